$wb = $excel.ActiveWorkbook

# --- General_MD (sheet3): insert 3 new rows (NOTEX_da/en/kl) above the
# "position" row, and append a new COPYRIGHT/YES row at the bottom. ---
$wsGeneral = $wb.Worksheets.Item("General_MD")

$wsGeneral.Rows.Item(10).Insert()
$wsGeneral.Rows.Item(10).Insert()
$wsGeneral.Rows.Item(10).Insert()

$wsGeneral.Range("A10").Value = "NOTEX_da"
$wsGeneral.Range("A11").Value = "NOTEX_en"
$wsGeneral.Range("A12").Value = "NOTEX_kl"

$wsGeneral.Range("B10").Value = "Tvungen fodnote"
$wsGeneral.Range("B10").WrapText = $true
$wsGeneral.Range("B11").Value = "Tvungen fodnote"
$wsGeneral.Range("B11").WrapText = $true
$wsGeneral.Range("B12").Value = "Tvungen fodnote"
$wsGeneral.Range("B12").WrapText = $true

# --- Variables_MD (sheet1): fill in note + domain values for BEXSTATEST row. ---
$wsVariables = $wb.Worksheets.Item("Variables_MD")
$wsVariables.Range("G3").Value = "dette er en fodnote til variabel"
$wsVariables.Range("H3").Value = "dette er en fodnote til variabel"
$wsVariables.Range("I3").Value = "dette er en fodnote til variabel"

# --- back to General_MD: append COPYRIGHT row at the new bottom. ---
$wsGeneral.Range("A41").Value = "COPYRIGHT"
$wsGeneral.Range("B41").Value = "YES"

# --- Variables_MD: domain columns. ---
$wsVariables.Range("M3").Value = "test"
$wsVariables.Range("N3").Value = "test2n"
$wsVariables.Range("O3").Value = "testkl"

# --- Fix up hyperlinks on General_MD: their anchors shifted down by the
# inserted rows (B35/B36/B37 -> B38/B39/B40), but the engine does not
# shift hyperlink anchors automatically, so rebuild them. ---
$wsGeneral.Hyperlinks.Delete()
$wsGeneral.Hyperlinks.Add($wsGeneral.Range("B38"), "http://www.stat.gl/bed202201/m1")
$wsGeneral.Hyperlinks.Add($wsGeneral.Range("B39"), "http://www.stat.gl/bee202201/m1")
$wsGeneral.Hyperlinks.Add($wsGeneral.Range("B40"), "http://www.stat.gl/ben202201/m1")
$wsGeneral.Range("B38").Style = "Hyperlink"
$wsGeneral.Range("B39").Style = "Hyperlink"
$wsGeneral.Range("B40").Style = "Hyperlink"

# --- Column widths on General_MD widened (closest achievable via ColumnWidth). ---
$wsGeneral.Columns.Item(1).ColumnWidth = 40.75
$wsGeneral.Columns.Item(2).ColumnWidth = 41.75

# --- Selections / active sheet: Variables_MD becomes the active tab
# (was General_MD), with new selections on each sheet. ---
$wsGeneral.Activate()
$wsGeneral.Range("B41").Select()

$wsVariables.Activate()
$wsVariables.Range("N3").Select()
